$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").ClearContents()
$ws.Range("F1").Value = "YoA"
$ws.Range("E2").ClearContents()

$ws.Range("E2").Select()
